$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value into a cell, forcing the destination to
# keep/borrow a given cell's visual style (font/alignment/wrap) while
# treating the new value strictly as text. This avoids two pitfalls of a
# naive "$ws.Range(...).Value = ..." assignment:
#   1) Excel's automatic literal-detection would silently convert a
#      date-shaped string like "01/01/2023" into a numeric date serial.
#   2) A brand-new, never-before-populated cell has no style of its own and
#      would otherwise fall back to the plain column style instead of
#      matching its row neighbours (the normal text style in column B, the
#      "changed/highlighted" red style in column C).
function Set-TextValue {
    param(
        [string]$CellAddr,
        [string]$StyleSourceAddr,
        [string]$Text
    )

    $target = $ws.Range($CellAddr)
    $source = $ws.Range($StyleSourceAddr)
    $scratch = $ws.Range("Z1")

    # Stage the text-forced value + number format on an out-of-the-way
    # scratch cell that starts from the desired style.
    $source.Copy() | Out-Null
    $scratch.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $scratch.NumberFormat = "@"
    $scratch.Value = $Text

    # Give the destination the full visual style first (font, alignment,
    # wrap text, etc.)...
    $source.Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null    # xlPasteFormats

    # ...then overlay just the value + number format from the scratch cell,
    # so the destination ends up with the exact same resolved style as the
    # source plus the text-typed value.
    $scratch.Copy() | Out-Null
    $target.PasteSpecial(-4163) | Out-Null    # xlPasteValuesAndNumberFormats

    $scratch.Clear() | Out-Null
}

# --- Ativação date: 01/01/2016 -> 01/01/2023. The same string previously
#     appeared in two places in the sheet (row 8 and, as leftover/stale
#     content, row 13); both must be updated since the shared text changes
#     for every cell that carried it. ---
Set-TextValue "B8" "B8" "01/01/2023"
Set-TextValue "C8" "C8" "01/01/2023"
Set-TextValue "B13" "B13" "01/01/2023"
Set-TextValue "C13" "C13" "01/01/2023"

# --- New English translations added alongside the existing "Objectives:",
#     "Short syllabus:" and "Syllabus:" rows (11, 14, 16). Column B borrows
#     the plain wrap-text style from column B's existing translated rows
#     (B10), column C borrows the red "changed" wrap-text style from C10. ---

$objectivesText = "Develop theoretical and practical knowledge of the manufacturing processes of equipment and devices required for the development of products and prototypes. Know the requirements and effects of manufacturing processes in order to allow, interact, create and execute projects throughout your professional life."
$shortSyllabusText = "Introduction to manufacturing processes. Material joining processes. Computer-aided design (CAD) review. Computer Aided Manufacturing (CAM). Flexible production systems. Rapid prototyping."
$syllabusText = "Classification of manufacturing processes. Foundry. Powder metallurgy. Machining: processes, fundamentals and economic conditions. Machine tools. Mechanical conformation. Material joining processes. Computer Aided Manufacturing (CAM). Programming languages for numerical control. Numerical control machine tools. Product manufacturing sequence. Notions of automation of manufacturing processes. Rapid prototyping. Rapid prototyping systems (solid, liquid and powder)."

Set-TextValue "B11" "B10" $objectivesText
Set-TextValue "C11" "C10" $objectivesText

Set-TextValue "B14" "B10" $shortSyllabusText
Set-TextValue "C14" "C10" $shortSyllabusText

Set-TextValue "B16" "B10" $syllabusText
Set-TextValue "C16" "C10" $syllabusText
